# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - mirror the style used by the other header cells (e.g. AC1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-47 get the team's season record: 91 wins, 71 losses, 0 ties
for ($row = 2; $row -le 47; $row++) {
    $ws.Cells.Item($row, 30).Value = 91
    $ws.Cells.Item($row, 31).Value = 71
    $ws.Cells.Item($row, 32).Value = 0
}
